$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns before column B (shifts old B -> E, carrying values & styles).
$ws.Range("B1:D1").EntireColumn.Insert()

# Match column widths of the inserted columns to column B/E (55.5 chars wide).
$ws.Columns.Item(2).ColumnWidth = 54.67
$ws.Columns.Item(3).ColumnWidth = 54.67
$ws.Columns.Item(4).ColumnWidth = 54.67

$ws.Cells.Item(1,"B").Value = "Jun_27"
$ws.Cells.Item(1,"C").Value = "Jun_26"
$ws.Cells.Item(1,"D").Value = "Jun_26"
$ws.Cells.Item(2,"B").Value = "UN"
$ws.Cells.Item(2,"C").Value = "UN"
$ws.Cells.Item(2,"D").Value = "UN"
$ws.Cells.Item(3,"B").Value = "UN"
$ws.Cells.Item(3,"C").Value = "UN"
$ws.Cells.Item(3,"D").Value = "UN"
$ws.Cells.Item(4,"B").Value = "UN"
$ws.Cells.Item(4,"C").Value = "UN"
$ws.Cells.Item(4,"D").Value = "UN"
$ws.Cells.Item(5,"B").Value = "UN"
$ws.Cells.Item(5,"C").Value = "UN"
$ws.Cells.Item(5,"D").Value = "UN"
$ws.Cells.Item(6,"B").Value = "UN"
$ws.Cells.Item(6,"C").Value = "UN"
$ws.Cells.Item(6,"D").Value = "UN"
$ws.Cells.Item(7,"B").Value = "UN"
$ws.Cells.Item(7,"C").Value = "UN"
$ws.Cells.Item(7,"D").Value = "UN"
$ws.Cells.Item(8,"B").Value = "UN"
$ws.Cells.Item(8,"C").Value = "UN"
$ws.Cells.Item(8,"D").Value = "UN"
$ws.Cells.Item(9,"B").Value = "UN"
$ws.Cells.Item(9,"C").Value = "UN"
$ws.Cells.Item(9,"D").Value = "UN"
$ws.Cells.Item(10,"B").Value = "UN"
$ws.Cells.Item(10,"C").Value = "UN"
$ws.Cells.Item(10,"D").Value = "UN"
$ws.Cells.Item(11,"B").Value = "UN"
$ws.Cells.Item(11,"C").Value = "UN"
$ws.Cells.Item(11,"D").Value = "UN"
$ws.Cells.Item(12,"B").Value = "UN"
$ws.Cells.Item(12,"C").Value = "UN"
$ws.Cells.Item(12,"D").Value = "UN"
$ws.Cells.Item(13,"B").Value = "UN"
$ws.Cells.Item(13,"C").Value = "UN"
$ws.Cells.Item(13,"D").Value = "UN"
$ws.Cells.Item(14,"B").Value = "UN"
$ws.Cells.Item(14,"C").Value = "UN"
$ws.Cells.Item(14,"D").Value = "UN"
$ws.Cells.Item(15,"B").Value = "UN"
$ws.Cells.Item(15,"C").Value = "UN"
$ws.Cells.Item(15,"D").Value = "UN"
$ws.Cells.Item(16,"B").Value = "UN"
$ws.Cells.Item(16,"C").Value = "UN"
$ws.Cells.Item(16,"D").Value = "UN"
$ws.Cells.Item(17,"B").Value = "UN"
$ws.Cells.Item(17,"C").Value = "UN"
$ws.Cells.Item(17,"D").Value = "UN"
$ws.Cells.Item(18,"B").Value = "UN"
$ws.Cells.Item(18,"C").Value = "UN"
$ws.Cells.Item(18,"D").Value = "UN"
$ws.Cells.Item(19,"B").Value = "UN"
$ws.Cells.Item(19,"C").Value = "UN"
$ws.Cells.Item(19,"D").Value = "UN"
$ws.Cells.Item(20,"B").Value = "UN"
$ws.Cells.Item(20,"C").Value = "UN"
$ws.Cells.Item(20,"D").Value = "UN"
$ws.Cells.Item(21,"B").Value = "UN"
$ws.Cells.Item(21,"C").Value = "6/19/2018,Reiterates,Buy,`$53.00"
$ws.Cells.Item(21,"D").Value = "6/19/2018,Reiterates,Buy,`$53.00"
$ws.Cells.Item(22,"B").Value = "UN"
$ws.Cells.Item(22,"C").Value = "UN"
$ws.Cells.Item(22,"D").Value = "UN"
$ws.Cells.Item(23,"B").Value = "UN"
$ws.Cells.Item(23,"C").Value = "UN"
$ws.Cells.Item(23,"D").Value = "UN"
$ws.Cells.Item(24,"B").Value = "UN"
$ws.Cells.Item(24,"C").Value = "UN"
$ws.Cells.Item(24,"D").Value = "UN"
$ws.Cells.Item(25,"B").Value = "UN"
$ws.Cells.Item(25,"C").Value = "UN"
$ws.Cells.Item(25,"D").Value = "UN"
$ws.Cells.Item(26,"B").Value = "UN"
$ws.Cells.Item(26,"C").Value = "UN"
$ws.Cells.Item(26,"D").Value = "UN"
$ws.Cells.Item(27,"B").Value = "UN"
$ws.Cells.Item(27,"C").Value = "UN"
$ws.Cells.Item(27,"D").Value = "UN"
$ws.Cells.Item(28,"A").Value = "Benchmark"
$ws.Cells.Item(28,"B").Value = "UN"
$ws.Cells.Item(28,"C").Value = "UN"
$ws.Cells.Item(28,"D").Value = "UN"
$ws.Cells.Item(29,"A").Value = "Evercore ISI"
$ws.Cells.Item(29,"B").Value = "UN"
$ws.Cells.Item(29,"C").Value = "UN"
$ws.Cells.Item(29,"D").Value = "UN"

Write-Output "done"
